$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 266.75
$ws.Range("I101").Value = 266.75
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 800.25
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = 821.75
$ws.Range("N101").ClearContents()
$ws.Range("H137").Value = 2243
$ws.Range("I137").Value = 1924.3334
$ws.Range("K137").Value = 5773.0002
$ws.Range("M137").Value = -3223.0002

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2585.4827
$ws.Range("J2").Value = 4502.5
$ws.Range("L2").Value = 4502.5
$ws.Range("N2").Value = -4728.5
$ws.Range("H61").Value = 34491790
$ws.Range("J61").Value = 47629140
$ws.Range("L61").Value = 47629140
$ws.Range("N61").Value = -47629564
$ws.Range("H69").Value = 250000
$ws.Range("J69").Value = 250000
$ws.Range("L69").Value = 250000
$ws.Range("N69").Value = -251498
$ws.Range("H72").Value = 250000
$ws.Range("J72").Value = 250000
$ws.Range("L72").Value = 750000
$ws.Range("N72").Value = -757488
$ws.Range("H74").Value = 47188.39
$ws.Range("I74").Value = 86110.25
$ws.Range("K74").Value = 86110.25
$ws.Range("M74").Value = -85236.25
$ws.Range("H77").Value = 47188.39
$ws.Range("I77").Value = 86110.25
$ws.Range("K77").Value = 430551.25
$ws.Range("M77").Value = -426183.25
$ws.Range("H116").Value = 2585.4827
$ws.Range("J116").Value = 4502.5
$ws.Range("L116").Value = 4502.5
$ws.Range("N116").Value = -9090.5
$ws.Range("H132").Value = 6532
$ws.Range("I132").Value = 3114.5334
$ws.Range("K132").Value = 9343.600199999999
$ws.Range("M132").Value = -6813.600199999999
$ws.Range("H136").Value = 34491790
$ws.Range("J136").Value = 47629140
$ws.Range("L136").Value = 142887420
$ws.Range("N136").Value = -142892520

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2585.4827
$ws.Range("J3").Value = 4502.5
$ws.Range("L3").Value = 4502.5
$ws.Range("N3").Value = -4730.5
$ws.Range("H86").Value = 96591.45
$ws.Range("I86").Value = 147214.42
$ws.Range("K86").Value = 147214.42
$ws.Range("M86").Value = -146091.42
$ws.Range("H89").Value = 96591.45
$ws.Range("I89").Value = 147214.42
$ws.Range("K89").Value = 736072.1000000001
$ws.Range("M89").Value = -730456.1000000001
$ws.Range("H134").Value = 5108193.5
$ws.Range("I134").Value = 10419774
$ws.Range("K134").Value = 31259322
$ws.Range("M134").Value = -31256787

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12429.031
$ws.Range("I31").Value = 5594.5
$ws.Range("J31").Value = 14707.208
$ws.Range("K31").Value = 5594.5
$ws.Range("L31").Value = 14707.208
$ws.Range("M31").Value = -5299.5
$ws.Range("N31").Value = -15297.208
$ws.Range("H34").Value = 12429.031
$ws.Range("I34").Value = 5594.5
$ws.Range("J34").Value = 14707.208
$ws.Range("K34").Value = 5594.5
$ws.Range("L34").Value = 14707.208
$ws.Range("M34").Value = -5392.5
$ws.Range("N34").Value = -15111.208
$ws.Range("H58").Value = 8147.185
$ws.Range("I58").Value = 5068.5713
$ws.Range("K58").Value = 5068.5713
$ws.Range("M58").Value = -4865.5713
$ws.Range("H132").Value = 9042.407
$ws.Range("I132").Value = 7427.5
$ws.Range("J132").Value = 10781.538
$ws.Range("K132").Value = 22282.5
$ws.Range("L132").Value = 32344.614
$ws.Range("M132").Value = -19752.5
$ws.Range("N132").Value = -37404.614
$ws.Range("H134").Value = 8215.392
$ws.Range("I134").Value = 3139.4
$ws.Range("K134").Value = 9418.2
$ws.Range("M134").Value = -6883.200000000001
$ws.Range("H136").Value = 8147.185
$ws.Range("I136").Value = 5068.5713
$ws.Range("K136").Value = 15205.7139
$ws.Range("M136").Value = -12655.7139
$ws.Range("H141").Value = 859097
$ws.Range("J141").Value = 1396828.4
$ws.Range("L141").Value = 1396828.4
$ws.Range("N141").Value = -1407188.4

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 125136.625
$ws.Range("I29").Value = 65.5
$ws.Range("J29").Value = 500350
$ws.Range("K29").Value = 196.5
$ws.Range("L29").Value = 1501050
$ws.Range("M29").Value = 80.5
$ws.Range("N29").Value = -1501604

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5228.5454
$ws.Range("I132").Value = 3161.318
$ws.Range("J132").Value = 9363
$ws.Range("K132").Value = 9483.954000000002
$ws.Range("L132").Value = 28089
$ws.Range("M132").Value = -6953.954000000002
$ws.Range("N132").Value = -33149

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6513.5713
$ws.Range("I40").Value = 5137.846
$ws.Range("K40").Value = 5137.846
$ws.Range("M40").Value = -5001.846
$ws.Range("H122").Value = 6257.974
$ws.Range("I122").Value = 5607.684
$ws.Range("J122").Value = 6875.75
$ws.Range("K122").Value = 16823.052
$ws.Range("L122").Value = 20627.25
$ws.Range("M122").Value = -14373.052
$ws.Range("N122").Value = -25527.25
$ws.Range("H132").Value = 19240358
$ws.Range("J132").Value = 14749.833
$ws.Range("L132").Value = 44249.499
$ws.Range("N132").Value = -49309.499
$ws.Range("H136").Value = 16243.412
$ws.Range("I136").Value = 9372.75
$ws.Range("K136").Value = 28118.25
$ws.Range("M136").Value = -25568.25

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 15010000
$ws.Range("I18").Value = 15010000
$ws.Range("K18").Value = 15010000
$ws.Range("M18").Value = -15009827
$ws.Range("H21").Value = 15000
$ws.Range("I21").Value = 15000
$ws.Range("K21").Value = 15000
$ws.Range("M21").Value = -14765
$ws.Range("H28").Value = 5750
$ws.Range("I28").Value = 5000
$ws.Range("K28").Value = 5000
$ws.Range("M28").Value = -4652
$ws.Range("H35").Value = 15000
$ws.Range("I35").Value = 15000
$ws.Range("K35").Value = 15000
$ws.Range("M35").Value = -14710
$ws.Range("H132").Value = 12770.424
$ws.Range("I132").Value = 10316.96
$ws.Range("J132").Value = 20437.5
$ws.Range("K132").Value = 30950.88
$ws.Range("L132").Value = 61312.5
$ws.Range("M132").Value = -28420.88
$ws.Range("N132").Value = -66372.5
$ws.Range("H136").Value = 24052234
$ws.Range("I136").Value = 47621084
$ws.Range("K136").Value = 142863252
$ws.Range("M136").Value = -142860702
